$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17: H17,J17,L17,N17
$ws.Range("H17").Value = 1432.2778
$ws.Range("J17").Value = 1966.6666
$ws.Range("L17").Value = 5899.9998
$ws.Range("N17").Value = -6235.9998
# Row 45: H45,J45,L45,N45
$ws.Range("H45").Value = 2414
$ws.Range("J45").Value = 2414
$ws.Range("L45").Value = 7242
$ws.Range("N45").Value = -7626
# Row 97: H97,J97,L97,N97
$ws.Range("H97").Value = 1529.1333
$ws.Range("J97").Value = 1529.1333
$ws.Range("L97").Value = 4587.3999
$ws.Range("N97").Value = -5579.3999
# Row 107: H107,I107,J107,K107,L107,M107,N107
$ws.Range("H107").Value = 532.7778
$ws.Range("I107").Value = 583.4
$ws.Range("J107").Value = 279.66666
$ws.Range("K107").Value = 583.4
$ws.Range("L107").Value = 279.66666
$ws.Range("M107").Value = 1336.6
$ws.Range("N107").Value = -4119.66666
# Row 112: H112,J112,L112,N112
$ws.Range("H112").Value = 2230.2917
$ws.Range("J112").Value = 2230.2917
$ws.Range("L112").Value = 6690.875100000001
$ws.Range("N112").Value = -8906.875100000001
# Row 125: H125,I125,K125,M125
$ws.Range("H125").Value = 2287.25
$ws.Range("I125").Value = 1849
$ws.Range("K125").Value = 16641
$ws.Range("M125").Value = -14181
# Row 137: H137,I137,J137,K137,L137,M137,N137
$ws.Range("H137").Value = 14928677
$ws.Range("I137").Value = 71431780
$ws.Range("J137").Value = 3330.1887
$ws.Range("K137").Value = 214295340
$ws.Range("L137").Value = 9990.5661
$ws.Range("M137").Value = -214292790
$ws.Range("N137").Value = -15090.5661
# Row 138: H138,I138,J138,K138,L138,M138,N138
$ws.Range("H138").Value = 4820.593
$ws.Range("I138").Value = 2577.7
$ws.Range("J138").Value = 5330.341
$ws.Range("K138").Value = 7733.099999999999
$ws.Range("L138").Value = 15991.023
$ws.Range("M138").Value = -2593.099999999999
$ws.Range("N138").Value = -26271.023
# Row 141: H141,I141,K141,M141
$ws.Range("H141").Value = 1982.75
$ws.Range("I141").Value = 977
$ws.Range("K141").Value = 2931
$ws.Range("M141").Value = 2249

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32: H32,I32,K32,M32
$ws.Range("H32").Value = 11081.325
$ws.Range("I32").Value = 9168.675999999999
$ws.Range("K32").Value = 9168.675999999999
$ws.Range("M32").Value = -8881.675999999999
# Row 61: H61,I61,J61,K61,L61,M61,N61
$ws.Range("H61").Value = 9792.214
$ws.Range("I61").Value = 11285
$ws.Range("J61").Value = 7801.8335
$ws.Range("K61").Value = 11285
$ws.Range("L61").Value = 7801.8335
$ws.Range("M61").Value = -11073
$ws.Range("N61").Value = -8225.833500000001
# Row 74: H74,I74,J74,K74,L74,M74,N74
$ws.Range("H74").Value = 4091.4
$ws.Range("I74").Value = 3313.2
$ws.Range("J74").Value = 4869.6
$ws.Range("K74").Value = 3313.2
$ws.Range("L74").Value = 4869.6
$ws.Range("M74").Value = -2439.2
$ws.Range("N74").Value = -6617.6
# Row 77: H77,I77,J77,K77,L77,M77,N77
$ws.Range("H77").Value = 4091.4
$ws.Range("I77").Value = 3313.2
$ws.Range("J77").Value = 4869.6
$ws.Range("K77").Value = 16566
$ws.Range("L77").Value = 24348
$ws.Range("M77").Value = -12198
$ws.Range("N77").Value = -33084
# Row 110: H110,I110,J110,K110,L110,M110,N110
$ws.Range("H110").Value = 4224.9585
$ws.Range("I110").Value = 3074.7778
$ws.Range("J110").Value = 7675.5
$ws.Range("K110").Value = 3074.7778
$ws.Range("L110").Value = 7675.5
$ws.Range("M110").Value = -1029.7778
$ws.Range("N110").Value = -11765.5
# Row 132: H132,I132,K132,M132
$ws.Range("H132").Value = 4408.6924
$ws.Range("I132").Value = 3800.5908
$ws.Range("K132").Value = 11401.7724
$ws.Range("M132").Value = -8871.7724
# Row 136: H136,I136,J136,K136,L136,M136,N136
$ws.Range("H136").Value = 9792.214
$ws.Range("I136").Value = 11285
$ws.Range("J136").Value = 7801.8335
$ws.Range("K136").Value = 33855
$ws.Range("L136").Value = 23405.5005
$ws.Range("M136").Value = -31305
$ws.Range("N136").Value = -28505.5005
# Row 137: H137,J137,L137,N137
$ws.Range("H137").Value = 67799.336
$ws.Range("J137").Value = 67799.336
$ws.Range("L137").Value = 67799.336
$ws.Range("N137").Value = -77999.336

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 58: H58,J58,L58,N58
$ws.Range("H58").Value = 19442.5
$ws.Range("J58").Value = 26390
$ws.Range("L58").Value = 26390
$ws.Range("N58").Value = -26978
# Row 59: H59,J59,L59,N59
$ws.Range("H59").Value = 76998.60000000001
$ws.Range("J59").Value = 89998.25
$ws.Range("L59").Value = 89998.25
$ws.Range("N59").Value = -91692.25
# Row 81: H81,J81,L81,N81
$ws.Range("H81").Value = 29889.25
$ws.Range("J81").Value = 29889.25
$ws.Range("L81").Value = 29889.25
$ws.Range("N81").Value = -32011.25
# Row 84: H84,J84,L84,N84
$ws.Range("H84").Value = 29889.25
$ws.Range("J84").Value = 29889.25
$ws.Range("L84").Value = 89667.75
$ws.Range("N84").Value = -100275.75
# Row 105: H105,I105,J105,K105,L105,M105,N105
$ws.Range("H105").Value = 14208.88
$ws.Range("I105").Value = 12226.167
$ws.Range("J105").Value = 19307.285
$ws.Range("K105").Value = 12226.167
$ws.Range("L105").Value = 19307.285
$ws.Range("M105").Value = -10479.167
$ws.Range("N105").Value = -22801.285
# Row 134: H134,I134,K134,M134
$ws.Range("H134").Value = 4465.0435
$ws.Range("I134").Value = 3981
$ws.Range("K134").Value = 11943
$ws.Range("M134").Value = -9408

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31: H31,I31,J31,K31,L31,M31,N31
$ws.Range("H31").Value = 37263.97
$ws.Range("I31").Value = 7927.8184
$ws.Range("J31").Value = 91046.914
$ws.Range("K31").Value = 7927.8184
$ws.Range("L31").Value = 91046.914
$ws.Range("M31").Value = -7632.8184
$ws.Range("N31").Value = -91636.914
# Row 34: H34,I34,J34,K34,L34,M34,N34
$ws.Range("H34").Value = 37263.97
$ws.Range("I34").Value = 7927.8184
$ws.Range("J34").Value = 91046.914
$ws.Range("K34").Value = 7927.8184
$ws.Range("L34").Value = 91046.914
$ws.Range("M34").Value = -7725.8184
$ws.Range("N34").Value = -91450.914
# Row 58: H58,I58,J58,K58,L58,M58,N58
$ws.Range("H58").Value = 6625.3335
$ws.Range("I58").Value = 1884.8
$ws.Range("J58").Value = 8995.6
$ws.Range("K58").Value = 1884.8
$ws.Range("L58").Value = 8995.6
$ws.Range("M58").Value = -1681.8
$ws.Range("N58").Value = -9401.6
# Row 134: H134,I134,K134,M134
$ws.Range("H134").Value = 4015.2144
$ws.Range("I134").Value = 2599.9167
$ws.Range("K134").Value = 7799.750100000001
$ws.Range("M134").Value = -5264.750100000001
# Row 136: H136,I136,J136,K136,L136,M136,N136
$ws.Range("H136").Value = 6625.3335
$ws.Range("I136").Value = 1884.8
$ws.Range("J136").Value = 8995.6
$ws.Range("K136").Value = 5654.4
$ws.Range("L136").Value = 26986.8
$ws.Range("M136").Value = -3104.4
$ws.Range("N136").Value = -32086.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 3: H3,I3,K3,M3
$ws.Range("H3").Value = 2426.2104
$ws.Range("I3").Value = 946.94116
$ws.Range("K3").Value = 2840.82348
$ws.Range("M3").Value = -2728.82348
# Row 18: H18,I18,J18,K18,L18,M18,N18
$ws.Range("H18").Value = 489.5
$ws.Range("I18").Value = 455.1111
$ws.Range("J18").Value = 799
$ws.Range("K18").Value = 1365.3333
$ws.Range("L18").Value = 2397
$ws.Range("M18").Value = -1196.3333
$ws.Range("N18").Value = -2735
# Row 56: H56,I56,K56,M56
$ws.Range("H56").Value = 7203.6
$ws.Range("I56").Value = 7203.6
$ws.Range("K56").Value = 7203.6
$ws.Range("M56").Value = -6673.6
# Row 109: H109,I109,K109,M109
$ws.Range("H109").Value = 2361.4
$ws.Range("I109").Value = 1734.8889
$ws.Range("K109").Value = 5204.6667
$ws.Range("M109").Value = -4164.6667
# Row 122: H122,I122,J122,K122,L122,M122,N122
$ws.Range("H122").Value = 2582.1
$ws.Range("I122").Value = 871.5454999999999
$ws.Range("J122").Value = 3572.4211
$ws.Range("K122").Value = 7843.9095
$ws.Range("L122").Value = 32151.7899
$ws.Range("M122").Value = -5393.9095
$ws.Range("N122").Value = -37051.7899
# Row 124: H124,I124,K124,M124
$ws.Range("H124").Value = 8010.8237
$ws.Range("I124").Value = 5933
$ws.Range("K124").Value = 17799
$ws.Range("M124").Value = -12889
# Row 132: H132,J132,L132,N132
$ws.Range("H132").Value = 5089.6816
$ws.Range("J132").Value = 5459.615
$ws.Range("L132").Value = 49136.535
$ws.Range("N132").Value = -54196.535
# Row 134: H134,I134,J134,K134,L134,M134,N134
$ws.Range("H134").Value = 3495.611
$ws.Range("I134").Value = 3113
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 9339
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -4269
$ws.Range("N134").Value = -40140

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132: H132,I132,K132,M132
$ws.Range("H132").Value = 4702.5713
$ws.Range("I132").Value = 2380.8
$ws.Range("K132").Value = 7142.400000000001
$ws.Range("M132").Value = -4612.400000000001
# Row 135: H135,J135,L135,N135
$ws.Range("H135").Value = 64728.617
$ws.Range("J135").Value = 64728.617
$ws.Range("L135").Value = 64728.617
$ws.Range("N135").Value = -74868.617

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 61: H61,I61,K61,M61
$ws.Range("H61").Value = 15286.857
$ws.Range("I61").Value = 12750.75
$ws.Range("K61").Value = 12750.75
$ws.Range("M61").Value = -12548.75
# Row 93: H93,I93,K93,M93
$ws.Range("H93").Value = 2550.3635
$ws.Range("I93").Value = 2550.3635
$ws.Range("K93").Value = 2550.3635
$ws.Range("M93").Value = -1302.3635
# Row 100: H100,I100,K100,M100
$ws.Range("H100").Value = 6310.3335
$ws.Range("I100").Value = 4855.5
$ws.Range("K100").Value = 4855.5
$ws.Range("M100").Value = -4314.5
# Row 113: H113,I113,K113,M113
$ws.Range("H113").Value = 15286.857
$ws.Range("I113").Value = 12750.75
$ws.Range("K113").Value = 12750.75
$ws.Range("M113").Value = -10580.75
# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 5447.5386
$ws.Range("I132").Value = 4692.9585
$ws.Range("J132").Value = 14502.5
$ws.Range("K132").Value = 14078.8755
$ws.Range("L132").Value = 43507.5
$ws.Range("M132").Value = -11548.8755
$ws.Range("N132").Value = -48567.5
# Row 136: H136,J136,L136,N136
$ws.Range("H136").Value = 6565.5
$ws.Range("J136").Value = 7129.2856
$ws.Range("L136").Value = 21387.8568
$ws.Range("N136").Value = -26487.8568

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 113: H113,I113,J113,K113,L113,M113,N113
$ws.Range("H113").Value = 294.07693
$ws.Range("I113").Value = 313.17392
$ws.Range("J113").Value = 147.66667
$ws.Range("K113").Value = 939.5217600000001
$ws.Range("L113").Value = 443.00001
$ws.Range("M113").Value = 1230.47824
$ws.Range("N113").Value = -4783.00001
# Row 126: H126,I126,K126,M126
$ws.Range("H126").Value = 3094.0625
$ws.Range("I126").Value = 1741.1111
$ws.Range("K126").Value = 5223.3333
$ws.Range("M126").Value = -2753.3333
# Row 132: H132,I132,J132,K132,L132,M132,N132
$ws.Range("H132").Value = 6460.125
$ws.Range("I132").Value = 3946
$ws.Range("J132").Value = 14002.5
$ws.Range("K132").Value = 11838
$ws.Range("L132").Value = 42007.5
$ws.Range("M132").Value = -9308
$ws.Range("N132").Value = -47067.5
# Row 136: H136,I136,J136,K136,L136,M136,N136
$ws.Range("H136").Value = 5271.143
$ws.Range("I136").Value = 2383.923
$ws.Range("J136").Value = 13612
$ws.Range("K136").Value = 7151.768999999999
$ws.Range("L136").Value = 40836
$ws.Range("M136").Value = -4601.768999999999
$ws.Range("N136").Value = -45936
